$d = $word.ActiveDocument

$replacements = @(
    @{old="340÷8="; new="133÷8="},
    @{old="687÷6="; new="367÷9="},
    @{old="178÷2="; new="685÷4="},
    @{old="638÷8="; new="370÷4="},
    @{old="832÷2="; new="826÷6="},
    @{old="764÷2="; new="970÷6="},
    @{old="411÷6="; new="488÷2="},
    @{old="742÷5="; new="661÷4="},
    @{old="982÷2="; new="526÷2="},
    @{old="235÷2="; new="555÷6="},
    @{old="710÷8="; new="733÷7="},
    @{old="603÷2="; new="867÷4="},
    @{old="756÷5="; new="914÷3="},
    @{old="523÷4="; new="362÷3="},
    @{old="986÷2="; new="122÷3="},
    @{old="301÷6="; new="594÷3="},
    @{old="573÷2="; new="451÷2="},
    @{old="944÷4="; new="657÷6="},
    @{old="925÷4="; new="123÷7="},
    @{old="882÷9="; new="781÷4="},
    @{old="994÷8="; new="224÷4="},
    @{old="968÷7="; new="971÷3="},
    @{old="298÷9="; new="879÷8="},
    @{old="241÷3="; new="691÷3="},
    @{old="282÷5="; new="930÷9="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
